# Inserts a new client row ("BETTER DRINKS") above the current row 17
# ("BRASA FOODIES"), pushing it (and every row below it) down by one,
# and updates the workbook-level "Clientes" defined name / dimension
# to cover the extra row (AE159 -> AE160).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Insert a brand-new blank row at position 17; existing rows 17..159
# shift down to 18..160 automatically (matches the diff exactly).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new client's data.
$ws.Range("A17").Value = 174
$ws.Range("B17").Value = "BETTER DRINKS"
$ws.Range("V17").Value = 44700
$ws.Range("W17").Value = 0
$ws.Range("X17").Value = 0
$ws.Range("Y17").Value = 0
$ws.Range("AA17").Value = 0
$ws.Range("AB17").Value = $false

# Update the "Clientes" defined name to reflect the new row count.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Clientes") {
        $n.RefersTo = "='Clientes'!`$A`$1:`$AE`$160"
    }
}
